$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-01")

# Rename the "Task Name" column header to "Name"
$ws.Range("C1").Value = "Name"

# Update the active cell selection (mirrors the saved sheet view state)
$ws.Range("C2").Select()

